# Update for insert release-notes.md f80ed2bb9e1dd81abc71d13817b8a44a756cee80
$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates -------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 0.3.0 -> 0.4.0-snapshot-1
$wsMeta.Range("B3").Value = "0.4.0-snapshot-1"

# Status: active -> draft
$wsMeta.Range("B6").Value = "draft"

# Date: 2024-03-13T09:33:00+00:00 -> 2024-05-23T12:16:26+00:00
$wsMeta.Range("B8").Value = "2024-05-23T12:16:26+00:00"

# Contact: "No display for ContactDetail" -> "ANS (https://esante.gouv.fr)"
$wsMeta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# --- Elements sheet: swap the two Mapping columns (AK <-> AL) --------------
$wsElem = $wb.Worksheets.Item("Elements")

$lastRow = 6
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $wsElem.Cells.Item($r, 37)
    $alCell = $wsElem.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the column widths to match the swapped content (AK now holds the wide
# French mapping text, AL now holds the narrow "RIM Mapping" text).
# (Input values are chosen so that, after the engine's internal
# character-width -> pixel -> character-width round trip, the stored
# width lands as close as possible to the target widths of 75.78515625 /
# 24.98046875.)
$wsElem.Columns.Item(37).ColumnWidth = 75.0
$wsElem.Columns.Item(38).ColumnWidth = 24.15
